# Updated cryptos list on Sat Apr 27 12:41:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.037.19'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.136.79'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.07'
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.19'
$ws.Range('E6').Value = '  -4.79%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.131.89'
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('E10').Value = '  -3.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.23'
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('E12').Value = '  -3.22%  '
$ws.Range('E13').Value = '  -4.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.10'
$ws.Range('E14').Value = '  -3.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.654.80'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.134.51'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.013.94'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.65'
$ws.Range('E19').Value = '  -3.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '470.08'
$ws.Range('E20').Value = '  -2.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.14'
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.696'
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.65'
$ws.Range('E23').Value = '  -0.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.93'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.94'
$ws.Range('E25').Value = '  -3.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.70'
$ws.Range('E27').Value = '  -2.53%  '
$ws.Range('E28').Value = '  -5.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.11'
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.92'
$ws.Range('E30').Value = '  -4.51%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.106'
$ws.Range('E33').Value = '  -6.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.52'
$ws.Range('E34').Value = '  -6.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.06'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '52.32'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0694'
$ws.Range('E38').Value = '  -9.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0385'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '415.70'
$ws.Range('E40').Value = '  -6.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.71'
$ws.Range('E41').Value = '  -9.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.18'
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.897.72'
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.111'
$ws.Range('E44').Value = '  -6.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.260'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  -6.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.34'
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.112'
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('E50').Value = '  -9.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.71'
$ws.Range('E51').Value = '  -0.80%  '
